# Weekly update: insert a new Puerro price observation for Vega Modelo de
# Temuco ahead of the existing row 201, shifting the rest of the table down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(201).Insert()

$ws.Cells.Item(201, 1).Value = 10
$ws.Cells.Item(201, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(201, 3).Value = "La Araucanía"
$ws.Cells.Item(201, 4).Value = 44746
$ws.Cells.Item(201, 5).Value = 9
$ws.Cells.Item(201, 6).Value = 100112005
$ws.Cells.Item(201, 7).Value = "Puerro"
$ws.Cells.Item(201, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(201, 9).Value = "Primera"
$ws.Cells.Item(201, 10).Value = 35
$ws.Cells.Item(201, 11).Value = 16000
$ws.Cells.Item(201, 12).Value = 16000
$ws.Cells.Item(201, 13).Value = 16000
$ws.Cells.Item(201, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(201, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(201, 16).Value = 1333
$ws.Cells.Item(201, 17).Value = 12
$ws.Cells.Item(201, 18).Value = "Hortaliza"
